$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'300.43"
$ws.Cells.Item(2, 5).Value = "'-0.52%"
$ws.Cells.Item(2, 7).Value = "'8"
$ws.Cells.Item(3, 4).Value = "'31.45"
$ws.Cells.Item(3, 5).Value = "'-1.17%"
$ws.Cells.Item(3, 7).Value = "'8"
$ws.Cells.Item(4, 4).Value = "'5.083"
$ws.Cells.Item(4, 5).Value = "'-1.93%"
$ws.Cells.Item(4, 7).Value = "'8"
$ws.Cells.Item(5, 4).Value = "'0.07857"
$ws.Cells.Item(5, 5).Value = "'-0.41%"
$ws.Cells.Item(5, 7).Value = "'8"
$ws.Cells.Item(6, 4).Value = "'2.261"
$ws.Cells.Item(6, 5).Value = "'-1.39%"
$ws.Cells.Item(6, 7).Value = "'8"
$ws.Cells.Item(7, 4).Value = "'7.803"
$ws.Cells.Item(7, 5).Value = "'-1.72%"
$ws.Cells.Item(7, 7).Value = "'8"
$ws.Cells.Item(8, 4).Value = "'3.836"
$ws.Cells.Item(8, 5).Value = "'-0.96%"
$ws.Cells.Item(8, 7).Value = "'8"
$ws.Cells.Item(9, 4).Value = "'0.9220"
$ws.Cells.Item(9, 5).Value = "'1.47%"
$ws.Cells.Item(9, 7).Value = "'8"
$ws.Cells.Item(10, 4).Value = "'0.1747"
$ws.Cells.Item(10, 5).Value = "'0.64%"
$ws.Cells.Item(10, 7).Value = "'8"
$ws.Cells.Item(11, 4).Value = "'0.07603"
$ws.Cells.Item(11, 5).Value = "'3.99%"
$ws.Cells.Item(11, 7).Value = "'8"
$ws.Cells.Item(12, 4).Value = "'0.09097"
$ws.Cells.Item(12, 5).Value = "'12.25%"
$ws.Cells.Item(12, 7).Value = "'8"
$ws.Cells.Item(13, 4).Value = "'0.03003"
$ws.Cells.Item(13, 5).Value = "'-3.35%"
$ws.Cells.Item(13, 7).Value = "'8"
$ws.Cells.Item(14, 5).Value = "'0.87%"
$ws.Cells.Item(14, 7).Value = "'8"
$ws.Cells.Item(15, 4).Value = "'0.001507"
$ws.Cells.Item(15, 5).Value = "'-0.86%"
$ws.Cells.Item(15, 7).Value = "'8"
$ws.Cells.Item(16, 4).Value = "'0.006013"
$ws.Cells.Item(16, 5).Value = "'-2.98%"
$ws.Cells.Item(16, 7).Value = "'8"
$ws.Cells.Item(17, 4).Value = "'3.475"
$ws.Cells.Item(17, 5).Value = "'-0.24%"
$ws.Cells.Item(17, 7).Value = "'8"
$ws.Cells.Item(18, 4).Value = "'2.243"
$ws.Cells.Item(18, 5).Value = "'0.04%"
$ws.Cells.Item(18, 7).Value = "'8"
$ws.Cells.Item(19, 5).Value = "'0.52%"
$ws.Cells.Item(19, 7).Value = "'8"
$ws.Cells.Item(20, 4).Value = "'0.1317"
$ws.Cells.Item(20, 5).Value = "'-2.01%"
$ws.Cells.Item(20, 7).Value = "'8"
$ws.Cells.Item(21, 4).Value = "'3.807"
$ws.Cells.Item(21, 5).Value = "'-18.76%"
$ws.Cells.Item(21, 7).Value = "'8"
$ws.Cells.Item(22, 4).Value = "'0.1709"
$ws.Cells.Item(22, 5).Value = "'2.40%"
$ws.Cells.Item(22, 7).Value = "'8"
$ws.Cells.Item(23, 4).Value = "'0.04615"
$ws.Cells.Item(23, 5).Value = "'-0.88%"
$ws.Cells.Item(23, 7).Value = "'8"
$ws.Cells.Item(24, 4).Value = "'0.001254"
$ws.Cells.Item(24, 5).Value = "'-1.91%"
$ws.Cells.Item(24, 7).Value = "'8"
$ws.Cells.Item(25, 4).Value = "'0.004470"
$ws.Cells.Item(25, 5).Value = "'0.14%"
$ws.Cells.Item(25, 7).Value = "'8"
$ws.Cells.Item(26, 4).Value = "'0.0001250"
$ws.Cells.Item(26, 5).Value = "'5.57%"
$ws.Cells.Item(26, 7).Value = "'8"
$ws.Cells.Item(27, 4).Value = "'0.0003387"
$ws.Cells.Item(27, 5).Value = "'-2.48%"
$ws.Cells.Item(27, 7).Value = "'8"
$ws.Cells.Item(28, 7).Value = "'8"
$ws.Cells.Item(29, 7).Value = "'8"
$ws.Cells.Item(30, 7).Value = "'8"
$ws.Cells.Item(31, 7).Value = "'8"
$ws.Cells.Item(32, 7).Value = "'8"
$ws.Cells.Item(33, 7).Value = "'8"
$ws.Cells.Item(34, 7).Value = "'8"
$ws.Cells.Item(35, 7).Value = "'8"
$ws.Cells.Item(36, 7).Value = "'8"
$ws.Cells.Item(37, 7).Value = "'8"
$ws.Cells.Item(38, 7).Value = "'8"
$ws.Cells.Item(39, 4).Value = "'0.01741"
$ws.Cells.Item(39, 5).Value = "'-7.13%"
$ws.Cells.Item(39, 7).Value = "'8"
$ws.Cells.Item(40, 4).Value = "'0.04632"
$ws.Cells.Item(40, 5).Value = "'0.70%"
$ws.Cells.Item(40, 7).Value = "'8"
$ws.Cells.Item(41, 4).Value = "'0.007056"
$ws.Cells.Item(41, 5).Value = "'-0.25%"
$ws.Cells.Item(41, 7).Value = "'8"
$ws.Cells.Item(42, 4).Value = "'0.1358"
$ws.Cells.Item(42, 5).Value = "'0.18%"
$ws.Cells.Item(42, 7).Value = "'8"
$ws.Cells.Item(43, 4).Value = "'0.002190"
$ws.Cells.Item(43, 5).Value = "'1.83%"
$ws.Cells.Item(43, 7).Value = "'8"
$ws.Cells.Item(44, 4).Value = "'0.009753"
$ws.Cells.Item(44, 5).Value = "'-9.14%"
$ws.Cells.Item(44, 7).Value = "'8"
$ws.Cells.Item(45, 4).Value = "'0.00006268"
$ws.Cells.Item(45, 5).Value = "'-3.29%"
$ws.Cells.Item(45, 7).Value = "'8"
$ws.Cells.Item(46, 5).Value = "'-1.31%"
$ws.Cells.Item(46, 7).Value = "'8"
$ws.Cells.Item(47, 5).Value = "'4.56%"
$ws.Cells.Item(47, 7).Value = "'8"
$ws.Cells.Item(48, 4).Value = "'0.7460"
$ws.Cells.Item(48, 5).Value = "'-9.08%"
$ws.Cells.Item(48, 7).Value = "'8"
$ws.Cells.Item(49, 4).Value = "'0.00002099"
$ws.Cells.Item(49, 5).Value = "'-1.31%"
$ws.Cells.Item(49, 7).Value = "'8"
$ws.Cells.Item(50, 4).Value = "'0.0001999"
$ws.Cells.Item(50, 5).Value = "'-1.31%"
$ws.Cells.Item(50, 7).Value = "'8"
$ws.Cells.Item(51, 7).Value = "'8"
